$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.14144504070282
$ws.Range("B1").Value = 2.219400644302368
$ws.Range("C1").Value = 10.84332275390625
$ws.Range("D1").Value = 2.360852241516113
$ws.Range("E1").Value = 1.275194048881531
